$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that existed before but should be empty after the edit
$ws.Range("A13").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental."
$ws.Range("C10").Value = "Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental."

$ws.Range("B13").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C13").Value = "1304060 - Maria das Graças de Almeida Felipe"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."
$ws.Range("C14").Value = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division."
$ws.Range("C15").Value = "Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."
$ws.Range("C16").Value = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope analysis of cells structure: optical and electron microscope.Structure and function of major organic molecules: carbohydrates, lipids, nucleic acids and proteins. Internal organization of the cell: membrane structure and function; intracelular compartments and protein sorting; vesicular traffic (endocytosis and exocytosis).Nucleus and genetic material organization: structure and functionCell cycle and cell division: mitosis and meiosisCell energy conversion: mitochondria and chloroplast."
$ws.Range("C17").Value = "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope analysis of cells structure: optical and electron microscope.Structure and function of major organic molecules: carbohydrates, lipids, nucleic acids and proteins. Internal organization of the cell: membrane structure and function; intracelular compartments and protein sorting; vesicular traffic (endocytosis and exocytosis).Nucleus and genetic material organization: structure and functionCell cycle and cell division: mitosis and meiosisCell energy conversion: mitochondria and chloroplast."

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "MF = média finalMF = (P1 + P2)/2"
$ws.Range("C20").Value = "MF = média finalMF = (P1 + P2)/2"

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. -De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006."
$ws.Range("C22").Value = "-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. -De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006."
